# Repeating experiment for regularized TPR from new master branch.
# Adds a new row (44) to the Experiments log describing EXP42 (run_id 39),
# a repeat of run_id 29 executed from the new master branch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy per-column formatting from the last "orange band" row (38) ---
#        so the new row 44 picks up the same fills/borders/wrap styles as
#        its neighbours, without touching columns that stay empty (G:J).
$ws.Range("A38:F38").Copy()
$ws.Range("A44:F44").PasteSpecial(-4122)
$ws.Range("K38:L38").Copy()
$ws.Range("K44:L44").PasteSpecial(-4122)

# --- 2. Cell values (order matters: it drives shared-string insertion order) ---

# A44: experiment description (rich text built up below)
$ws.Range("A44").Value = "Just TPR no LSTM in `r`nphrase embedding layer `r`nbatchsize = 40. With visualizations. With regularization. Regularization weights=0.00001 [Repeating experiment run_id 29 from new master branch for repeatability purposes, running from QA_TPR_for_Run]. "

# D44: logfile name
$ws.Range("D44").Value = "EXP42.txt"

# B44: command to run
$ws.Range("B44").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --justTPR True --TPRregularizer1 True --TPRvis True --cF 0.00001 --cR 0.00001 --batch_size 40 --run_id 39 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP42.txt"

# C44: Machine / GPU (reuses existing shared string "DLT1 / 2")
$ws.Range("C44").Value = "DLT1 / 2"

# E44: run_id, F44: pane number in tmux
$ws.Range("E44").Value = 39
$ws.Range("F44").Value = 2

# L44: date
$ws.Range("L44").Value = [DateTime]"2017-03-21"
$ws.Range("L44").NumberFormat = "d-mmm-yy"
$ws.Range("L44").Interior.Color = 49407

# --- 3. Rich-text formatting for A44: alternating bold "labels" and plain text ---
$a44 = $ws.Range("A44")

# "With visualizations"
$a44.Characters(64, 19).Font.Bold = $true
$a44.Characters(64, 19).Font.ColorIndex = -4105

# ". "
$a44.Characters(83, 2).Font.Bold = $false
$a44.Characters(83, 2).Font.Size = 11
$a44.Characters(83, 2).Font.Name = "Calibri"
$a44.Characters(83, 2).Font.ColorIndex = -4105

# "With regularization"
$a44.Characters(85, 19).Font.Bold = $true
$a44.Characters(85, 19).Font.ColorIndex = -4105

# ". "
$a44.Characters(104, 2).Font.Bold = $false
$a44.Characters(104, 2).Font.Size = 11
$a44.Characters(104, 2).Font.Name = "Calibri"
$a44.Characters(104, 2).Font.ColorIndex = -4105

# "Regularization weights=0.00001"
$a44.Characters(106, 30).Font.Bold = $true
$a44.Characters(106, 30).Font.ColorIndex = -4105

# " [Repeating experiment run_id 29 from new master branch for repeatability purposes, running from QA_TPR_for_Run]. "
$a44.Characters(136, 114).Font.Bold = $false
$a44.Characters(136, 114).Font.Size = 11
$a44.Characters(136, 114).Font.Name = "Calibri"
$a44.Characters(136, 114).Font.ColorIndex = -4105

# --- 4. Row height + selection/view state ---
$ws.Rows("44:44").RowHeight = 195

$ws.Range("A44").Select()
